# "Asset List update & Uren Registratie"
#
# Uren Registratie Game-Lab sheet — Week 14 (rows 73-79) attendance-hours
# grid update: a few students who were previously marked 0 hours on
# Wednesday/Thursday now have hours logged, and the manually-entered
# "Totaal lesuren" (B79) for that week goes from 8 to 12. Everything else
# (K2, L2:L8, M2:M10, C79:I79, ...) is formula-driven and recalculates
# automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Woensdag (Wednesday) — row 76: Zinedine, Robin, Sam now logged 2 hours.
$ws.Range("D76").Value = 2   # Zinedine
$ws.Range("E76").Value = 2   # Robin
$ws.Range("H76").Value = 2   # Sam

# Donderdag (Thursday) — row 77: Rief, Zinedine, Robin, Marc logged 2 hours,
# Sam logged 1 hour.
$ws.Range("C77").Value = 2   # Rief
$ws.Range("D77").Value = 2   # Zinedine
$ws.Range("E77").Value = 2   # Robin
$ws.Range("G77").Value = 2   # Marc
$ws.Range("H77").Value = 1   # Sam

# Manually-entered "Totaal lesuren" for Week 14 increases from 8 to 12.
$ws.Range("B79").Value = 12

# Scroll the view down a bit and leave the cursor on the newly-edited block.
$win = $excel.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 1
[void]$ws.Range("K76").Select()

[void]$excel.Calculate()
